$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 12 data (new issue entry) - copy formats from the row above (style s="1")
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A12").Value = "Post-processing"
$ws.Range("B12").Value = "High"
$ws.Range("C12").Value = "Update expansion processing tracking to have a column each for rcmcs_processed, thermo_processed"
$ws.Range("D12").Value = "Change logic of process all to accommodate this"
$ws.Range("A12:D12").RowHeight = 45

# Add new Status column - copy header/body styles from column A of the corresponding rows
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "Status"

$ws.Range("A2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F2").Value = "Posted to GH"

$ws.Range("A3").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F3").Value = "Posted to GH"

# Widen column F
$ws.Columns.Item(6).ColumnWidth = 21

# Re-apply autofilter over the new range (this clears old sort state attached to filter)
$ws.AutoFilterMode = $false
$ws.Range("A1:F12").AutoFilter() | Out-Null

# Keep the _FilterDatabase defined name in sync with the new filter range
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$F`$12"

# Update the selection to match the recorded state
$ws.Range("A1:F1048576").Select() | Out-Null

$wb.Save()
